# "multiple shift types per shift"
# - Add a new trailing shift-type row (row 8) to the "Shifts" sheet,
#   continuing on from the existing row 7 block with its own start/end
#   time, shift type ("Standard,Reference") and label ("Project").
# - Widen column D on the Shifts sheet to fit the new, longer shift-type
#   text.
# - Make the Shifts sheet the active tab/selection (instead of Staff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Shifts"

# --- New row 8 values -----------------------------------------------
$ws.Range("B8").Value = 0.458333333333333   # 11:00 start
$ws.Range("C8").Value = 0.5                 # 12:00 end
$ws.Range("D8").Value = "Standard,Reference"
$ws.Range("F8").Value = "Project"

# Match D8's formatting/style to the other data rows in column D
# (style index 10 rather than the bare column-default style).
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats

# --- Widen column D to fit the longer shift-type text -----------------
$ws.Columns.Item(4).ColumnWidth = 25.15

# --- Make "Shifts" the active sheet/selection --------------------------
$ws.Range("C5").Select() | Out-Null
